# Clean-up of input tables
#
# The "efficiency class" lookup table had a block of rows (building
# component 2, option 21) whose "value" column used an irregular step
# sequence left over from an earlier draft. This corrects that column to
# the same clean progression used elsewhere in the table, and normalizes
# the header row styling (explicit black font) that the rest of the
# clean-up pass applied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected "value" column for id_building_component=2 / option=21 ---
$fixedValues = @{
    17 = 1.2
    18 = 1.4
    19 = 1.6
    20 = 1.8
    21 = 2
    22 = 2.5
    23 = 3
    24 = 4
    25 = 5
}

foreach ($row in $fixedValues.Keys) {
    $ws.Cells.Item($row, 5).Value = $fixedValues[$row]
}

# --- Normalize header row font to explicit black (matches clean-up pass) ---
$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Color = 0
